$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 24, shifting the existing data (rows 24-40) down
# to rows 25-41.
$ws.Rows(24).Insert()

# Populate the newly-inserted row 24 with a new weekly record for
# "Orégano" at "Vega Central Mapocho de Santiago".
$ws.Cells.Item(24, 1).Value = 9
$ws.Cells.Item(24, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(24, 3).Value = "Metropolitana"
$ws.Cells.Item(24, 4).Value = 44539
$ws.Cells.Item(24, 5).Value = 13
$ws.Cells.Item(24, 6).Value = 100112029
$ws.Cells.Item(24, 7).Value = "Orégano"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 16
$ws.Cells.Item(24, 11).Value = 9000
$ws.Cells.Item(24, 12).Value = 10000
$ws.Cells.Item(24, 13).Value = 9500
$ws.Cells.Item(24, 14).Value = "`$/docena de atados"
$ws.Cells.Item(24, 15).Value = "Región Metropolitana"
$ws.Cells.Item(24, 16).Value = 3167
$ws.Cells.Item(24, 17).Value = 3
$ws.Cells.Item(24, 18).Value = "Hortaliza"

# Row 35 (Volumen=160) shifted down to row 36; the stray extra digit in
# Volumen gets corrected to 16 at its new location.
$ws.Cells.Item(36, 10).Value = 16
